$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Red Zone" column (Q) ---
$ws.Range("Q1").Value = "Red Zone"

$redZoneByRow = @{
    2=1; 3=1; 4=1; 5=1
    6=0; 7=0; 8=0; 9=0
    10=0; 11=0; 12=0; 13=0
    14=1; 15=1; 16=1; 17=1
    18=1; 19=1; 20=1; 21=1
    22=1; 23=1; 24=1; 25=1
    26=1; 27=1; 28=1; 29=1
    30=1; 31=1; 32=1; 33=1
    34=0; 35=0; 36=0; 37=0
    38=1; 39=1; 40=1; 41=1
    42=0; 43=0; 44=0; 45=0
    46=0; 47=0; 48=0; 49=0
}

foreach ($r in $redZoneByRow.Keys) {
    $ws.Cells.Item($r, 17).Value = $redZoneByRow[$r]
}

# --- Update the view: zoom to 70%, select O18, drop the old frozen topLeftCell/selection ---
$ws.Range("O18").Select() | Out-Null
$excel.ActiveWindow.Zoom = 70

# --- Page setup: portrait orientation (adds a pageSetup element to the sheet) ---
$ws.PageSetup.Orientation = 1
